# Week 8 testing spreadsheet update
# - adds model-probability text for games in rows 2,3,5,6,7 (cols E:Q)
# - flips the "correct pick" indicator (col B, and D5) to 1 for those rows
# - re-activates "Season results" as the visible tab/selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 8 results")

# --- Numeric cell updates (B/D columns) ---
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1

# --- Clear explicit style on A5 / A7 (was "Neutral"-ish s="3", becomes default Normal) ---
$ws.Range("A5").Style = "Normal"
$ws.Range("A7").Style = "Normal"

# --- Row 3: fill columns E:Q with the ON001/ON002 model-probability text ---
$ws.Range("E3").Value = 'ON001(home win): 0.0003 / ON002(away win): 0.9996'
$ws.Range("F3").Value = 'ON001(home win): 0.9905  / ON002(away win): 0.0094'
$ws.Range("G3").Value = 'ON001(home win): 0.0023  / ON002(away win): 0.9976'
$ws.Range("H3").Value = 'ON001(home win): 0.1874  / ON002(away win): 0.8125'
$ws.Range("I3").Value = 'ON001(home win): 0.9901  / ON002(away win): 0.0098'
$ws.Range("J3").Value = 'ON001(home win): 0.9390  / ON002(away win): 0.0609'
$ws.Range("K3").Value = 'ON001(home win): 0.9996  / ON002(away win): 0.0003'
$ws.Range("L3").Value = 'ON001(home win): 0.1096  / ON002(away win): 0.8903'
$ws.Range("M3").Value = 'ON001(home win): 0.9918  / ON002(away win): 0.0081'
$ws.Range("N3").Value = 'ON001(home win): 0.9986  / ON002(away win): 0.0013'
$ws.Range("O3").Value = 'ON001(home win): 0.2467  / ON002(away win): 0.7532'
$ws.Range("P3").Value = 'ON001(home win): 0.0002  / ON002(away win): 0.9997'
$ws.Range("Q3").Value = 'ON001(home win): 0.0001  / ON002(away win): 0.9998'
$ws.Range("E3").Style = "Good"

# --- Row 2: fill columns E:Q with the ON001/ON002 model-probability text ---
$ws.Range("E2").Value = 'ON001(home win): 1.9e-9  / ON002(away win): 0.9999'
$ws.Range("F2").Value = 'ON001(home win): 0.9999  / ON002(away win): 0.00003'
$ws.Range("G2").Value = 'ON001(home win): 4.1e-9  / ON002(away win): 0.9999'
$ws.Range("H2").Value = 'ON001(home win): 0.0225  / ON002(away win): 0.9771'
$ws.Range("I2").Value = 'ON001(home win): 0.9987  / ON002(away win): 0.0012'
$ws.Range("J2").Value = 'ON001(home win): 0. 0727 / ON002(away win): 0.9234'
$ws.Range("K2").Value = 'ON001(home win): 0.9999  / ON002(away win): 1.08e-8'
$ws.Range("L2").Value = 'ON001(home win): 0.8273  / ON002(away win): 0.1689'
$ws.Range("M2").Value = 'ON001(home win): 0.9994  / ON002(away win): 0.0004'
$ws.Range("N2").Value = 'ON001(home win): 0.9999  / ON002(away win): 2.0e-8'
$ws.Range("O2").Value = 'ON001(home win): 0.0002  / ON002(away win): 0.9997'
$ws.Range("P2").Value = 'ON001(home win): 0.0001  / ON002(away win): 0.9998'
$ws.Range("Q2").Value = 'ON001(home win): 0.0012  / ON002(away win): 0.9987'
$ws.Range("E2").Style = "Good"

# --- Row 5: fill columns E:Q with the ON001/ON002 model-probability text ---
$ws.Range("E5").Value = 'ON001(home win): 2.3e-8  / ON002(away win): 0.9999'
$ws.Range("F5").Value = 'ON001(home win): 0.9999  / ON002(away win): 0.0000003'
$ws.Range("G5").Value = 'ON001(home win): 1.8e-11  / ON002(away win): 0.9999'
$ws.Range("H5").Value = 'ON001(home win): 0.4729  / ON002(away win): 0.5269'
$ws.Range("I5").Value = 'ON001(home win): 0.9905  / ON002(away win): 0.0091'
$ws.Range("J5").Value = 'ON001(home win): 0.6946  / ON002(away win): 0.3089'
$ws.Range("K5").Value = 'ON001(home win): 0.9999  / ON002(away win): 0.00007'
$ws.Range("L5").Value = 'ON001(home win): 0.0063  / ON002(away win): 0.9939'
$ws.Range("M5").Value = 'ON001(home win): 0.6489  / ON002(away win): 0.3507'
$ws.Range("N5").Value = 'ON001(home win): 0.9998  / ON002(away win): 0.0001'
$ws.Range("O5").Value = 'ON001(home win): 0.0010  / ON002(away win): 0.9990'
$ws.Range("P5").Value = 'ON001(home win): 1.1e-9  / ON002(away win): 0.9999'
$ws.Range("Q5").Value = 'ON001(home win): 1.2e-11  / ON002(away win): 0.9999'
$ws.Range("E5").Style = "Good"

# --- Row 6: fill columns E:Q with the ON001/ON002 model-probability text ---
$ws.Range("E6").Value = 'ON001(home win): 3.8e-9  / ON002(away win): 0.9999'
$ws.Range("F6").Value = 'ON001(home win): 0.9997  / ON002(away win): 0.0002'
$ws.Range("G6").Value = 'ON001(home win): 1.7e-8  / ON002(away win): 0.9999'
$ws.Range("H6").Value = 'ON001(home win): 0.2743  / ON002(away win): 0.7331'
$ws.Range("I6").Value = 'ON001(home win): 0.9999  / ON002(away win): 0.00002'
$ws.Range("J6").Value = 'ON001(home win): 0.9162  / ON002(away win): 0.0832'
$ws.Range("K6").Value = 'ON001(home win): 0.9999  / ON002(away win): 0.000002'
$ws.Range("L6").Value = 'ON001(home win): 0.0001  / ON002(away win): 0.9998'
$ws.Range("M6").Value = 'ON001(home win): 0.9932  / ON002(away win): 0.0068'
$ws.Range("N6").Value = 'ON001(home win): 0.9945  / ON002(away win): 0.0052'
$ws.Range("O6").Value = 'ON001(home win): 0.00002  / ON002(away win): 0.9999'
$ws.Range("P6").Value = 'ON001(home win): 9.9e-7  / ON002(away win): 0.9999'
$ws.Range("Q6").Value = 'ON001(home win): 6.1e-7  / ON002(away win): 0.9999'
$ws.Range("E6").Style = "Good"

# --- Row 7: fill columns E:Q with the ON001/ON002 model-probability text ---
$ws.Range("E7").Value = 'ON001(home win): 6.7e-14  / ON002(away win): 0.9999'
$ws.Range("F7").Value = 'ON001(home win): 0.9999  / ON002(away win): 3.6e-10'
$ws.Range("G7").Value = 'ON001(home win): 7.3e-14  / ON002(away win): 0.999'
$ws.Range("H7").Value = 'ON001(home win): 0.000002  / ON002(away win): 0.9999'
$ws.Range("I7").Value = 'ON001(home win): 0.9999  / ON002(away win):1.1e-11'
$ws.Range("J7").Value = 'ON001(home win): 0.1677  / ON002(away win): 0.8322'
$ws.Range("K7").Value = 'ON001(home win): 1.0  / ON002(away win): 5.5e-18'
$ws.Range("L7").Value = 'ON001(home win): 0.9999  / ON002(away win): 1.9e-9'
$ws.Range("M7").Value = 'ON001(home win): 0.9999  / ON002(away win): 1.9e-11'
$ws.Range("N7").Value = 'ON001(home win): 1.0  / ON002(away win): 1.2e-25'
$ws.Range("O7").Value = 'ON001(home win): 0.0350  / ON002(away win): 0.9649'
$ws.Range("P7").Value = 'ON001(home win): 6.8e-14  / ON002(away win): 0.9999'
$ws.Range("Q7").Value = 'ON001(home win): 1.5e-13  / ON002(away win): 0.9999'
$ws.Range("E7").Style = "Good"

# --- H5 is the one outlier cell styled "Neutral" instead of "Good" ---
$ws.Range("H5").Style = "Neutral"

# --- Sheet view / selection: Week 8 sheet keeps selection on Q7 but is no longer the active tab ---
[void]$ws.Range("Q7").Select()

$wsSeason = $wb.Worksheets.Item("Season results")
[void]$wsSeason.Activate()
[void]$wsSeason.Range("E4").Select()
